$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The monthly mailing list refresh: all "201711_..._FSCBook.pdf" attachment
# file names roll forward to "201712_..._FSCBook.pdf". These live in column L
# (rows 3-19); row 2's L/G cells and every row's G cell are driven by
# formulas referencing column L, so they recalculate automatically.
$codes = @(
    "ATLTC",
    "CALGY",
    "EDMON",
    "LONDN",
    "MEDIC",
    "MNTRL",
    "NWFLD",
    "OTTWA",
    "PRAIR",
    "QUEBC",
    "TORNT1",
    "TORNT2",
    "TORNT3",
    "TORNT4",
    "TORNT5",
    "VACVR",
    "VICTR"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 3
    $ws.Range("L$row").Value = "201712_" + $codes[$i] + "_FSCBook.pdf"
}

# Reflect the new active cell/selection recorded in the sheet view.
$ws.Range("A2").Select()
